$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts Late / heading / Outstanding
# one column to the right: N->O, O->P, P->Q) for the Loan RBI / Variable
# Instalments layout change.
$ws.Columns("N").Insert()

# Match the width of the neighbouring "In Advance" column (M) as closely as
# this host's ColumnWidth setter allows.
$ws.Columns("N").ColumnWidth = 9.83

# Activate the "Repayment schedule" sheet/tab and move the selection to R11,
# matching the saved view state.
$ws.Activate() | Out-Null
$ws.Range("R11").Select() | Out-Null
